# fix(gui) step 1 and 2
# - bump the sheet date by one day
# - update the unit prices in the "Bisagras FICHAS" price list

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: date in A1 advances by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: updated prices for B-522 / B-525 / B-526 / B-529
$ws.Range("D33").Value = 445.44
$ws.Range("D34").Value = 487.2
$ws.Range("D35").Value = 546.36
$ws.Range("D36").Value = 664.6799999999999
